$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.123485088348389
$ws.Range("B1").Value = 3.476093769073486
$ws.Range("C1").Value = 4.765398979187012
$ws.Range("D1").Value = 2.173427581787109
$ws.Range("E1").Value = 1.346295356750488
